# REAS_scaling_mapping.xlsx update
# Commit: General updates to include BC and OC scaling (Japan / REAS inventories),
# adding two new country/sector comment rows to the "year" sheet.

$wb = $excel.ActiveWorkbook

$wsMap  = $wb.Worksheets.Item("map")
$wsYear = $wb.Worksheets.Item("year")

# --- "year" sheet: add two new data rows (3 and 4) -------------------------
# Row 3: Afghanistan / road sector, big jump in road emissions after 2008.
# Write the comment text first so new shared-string entries are appended in
# the same order the target workbook uses them (Comment, iso, iso, Comment).
$wsYear.Range("H3").Value = "Big jump in road emissions after 2008"
$wsYear.Range("A3").Value = "afg"
$wsYear.Range("B3").Value = "road"
$wsYear.Range("C3").Value = 1960
$wsYear.Range("D3").Value = 2008
$wsYear.Range("E3").Value = "NA"
$wsYear.Range("F3").Value = 2000
$wsYear.Range("G3").Value = 2008

# Row 4: Cambodia / RCO sector, doubling of residential emissions 2008->2009.
$wsYear.Range("A4").Value = "khm"
$wsYear.Range("B4").Value = "RCO"
$wsYear.Range("C4").Value = 1960
$wsYear.Range("D4").Value = 2008
$wsYear.Range("E4").Value = "NA"
$wsYear.Range("F4").Value = 2000
$wsYear.Range("G4").Value = 2008
$wsYear.Range("H4").Value = "Doubling of residential emissions between 2008 and 2009"

# --- Sheet view / selection bookkeeping ------------------------------------
# Update "map" sheet's frozen-pane selection (scroll position itself is not
# exposed through the object model, only the active selection/sqref).
[void]$wsMap.Select()
[void]$wsMap.Range("A29:XFD35").Select()

# Restore "year" as the active sheet/tab and set its new selection, matching
# the workbook's saved state (year tab selected, cell E22 highlighted).
[void]$wsYear.Select()
[void]$wsYear.Range("E22").Select()
